# chore: update Sheets via scheduled runner
# Refresh computed profit columns (H-N) on the Odin_Profits workbook
# for the rows whose market data changed since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 533
$ws.Range("I5").Value = 599.5
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 599.5
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = -487.5
$ws.Range("N5").Value = -624

$ws.Range("H80").Value = 21530
$ws.Range("I80").Value = 12250
$ws.Range("J80").Value = 26170
$ws.Range("K80").Value = 12250
$ws.Range("L80").Value = 26170
$ws.Range("M80").Value = -11252
$ws.Range("N80").Value = -28166

$ws.Range("H83").Value = 21530
$ws.Range("I83").Value = 12250
$ws.Range("J83").Value = 26170
$ws.Range("K83").Value = 36750
$ws.Range("L83").Value = 78510
$ws.Range("M83").Value = -31758
$ws.Range("N83").Value = -88494

$ws.Range("H117").Value = 17624
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 17624
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 17624
$ws.Range("N117").Value = -26802

$ws.Range("H118").Value = 399999
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 399999
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 399999
$ws.Range("N118").Value = -403313

$ws.Range("H119").Value = 20674.25
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 20674.25
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 20674.25
$ws.Range("N119").Value = -30350.25

$ws.Range("H131").Value = 103198.8
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 103198.8
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 103198.8
$ws.Range("N131").Value = -113278.8

$ws.Range("H133").Value = 176482.83
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 176482.83
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 176482.83
$ws.Range("N133").Value = -181542.83

$ws.Range("H138").Value = 40429
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 40429
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 40429
$ws.Range("N138").Value = -50709

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 533
$ws.Range("I4").Value = 599.5
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 599.5
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -484.5
$ws.Range("N4").Value = -630

$ws.Range("H134").Value = 1620552.4
$ws.Range("I134").Value = 1673654.1
$ws.Range("J134").Value = 27500
$ws.Range("K134").Value = 5020962.300000001
$ws.Range("L134").Value = 82500
$ws.Range("M134").Value = -5018427.300000001
$ws.Range("N134").Value = -87570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13516049
$ws.Range("I16").Value = 29413518
$ws.Range("J16").Value = 3200.5
$ws.Range("K16").Value = 29413518
$ws.Range("L16").Value = 3200.5
$ws.Range("M16").Value = -29413231
$ws.Range("N16").Value = -3774.5

$ws.Range("H113").Value = 13516049
$ws.Range("I113").Value = 29413518
$ws.Range("J113").Value = 3200.5
$ws.Range("K113").Value = 29413518
$ws.Range("L113").Value = 3200.5
$ws.Range("M113").Value = -29411348
$ws.Range("N113").Value = -7540.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 91.59999999999999
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 91.59999999999999
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 91.59999999999999
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -317.6

$ws.Range("H122").Value = 4743.5757
$ws.Range("I122").Value = 3277.7144
$ws.Range("J122").Value = 7308.8335
$ws.Range("K122").Value = 9833.143199999999
$ws.Range("L122").Value = 21926.5005
$ws.Range("M122").Value = -7383.143199999999
$ws.Range("N122").Value = -26826.5005

$ws.Range("H138").Value = 80000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 80000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

$ws.Range("H139").Value = 188162.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 188162.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 188162.5
$ws.Range("N139").Value = -198442.5

$ws.Range("H141").Value = 51659.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 51659.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 51659.25
$ws.Range("N141").Value = -62019.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 45000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 45000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 45000
$ws.Range("N36").Value = -46124

$ws.Range("H98").Value = 69177
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 69177
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 69177
$ws.Range("N98").Value = -75167

$ws.Range("H127").Value = 25897.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 25897.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 25897.5
$ws.Range("N127").Value = -35817.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 200429
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 200429
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 200429
$ws.Range("N46").Value = -200891

$ws.Range("H104").Value = 28060
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 28060
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 28060
$ws.Range("N104").Value = -35048

$ws.Range("H112").Value = 49500
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 49500
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 49500
$ws.Range("N112").Value = -52454

$ws.Range("H116").Value = 55000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 55000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 55000
$ws.Range("N116").Value = -64178

$ws.Range("H118").Value = 60000
$ws.Range("I118").Value = 40000
$ws.Range("J118").Value = 80000
$ws.Range("K118").Value = 40000
$ws.Range("L118").Value = 80000
$ws.Range("M118").Value = -38343
$ws.Range("N118").Value = -83314

$ws.Range("H121").Value = 9258.333000000001
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 9258.333000000001
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 9258.333000000001
$ws.Range("N121").Value = -12752.333

$ws.Range("H134").Value = 200429
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 200429
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 601287
$ws.Range("N134").Value = -606357

$ws.Range("H135").Value = 84444
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 84444
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 84444
$ws.Range("N135").Value = -94584

$ws.Range("H136").Value = 8780497
$ws.Range("I136").Value = 11911157
$ws.Range("J136").Value = 14650.733
$ws.Range("K136").Value = 35733471
$ws.Range("L136").Value = 43952.199
$ws.Range("M136").Value = -35730921
$ws.Range("N136").Value = -49052.199

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 98515.42999999999
$ws.Range("I141").Value = 79995
$ws.Range("J141").Value = 101602.164
$ws.Range("K141").Value = 79995
$ws.Range("L141").Value = 101602.164
$ws.Range("M141").Value = -74815
$ws.Range("N141").Value = -123085
